$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'29.451.15"
$ws.Range("E2").Value = "  +0.71%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.873.95"
$ws.Range("E3").Value = "  +0.95%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.11%  "

# Row 5 - XRP
$ws.Range("D5").Value = "'0.7140"
$ws.Range("E5").Value = "  +1.67%  "

# Row 6 - BNB
$ws.Range("E6").Value = "  +1.79%  "

# Row 7
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("D8").Value = "'0.07903"
$ws.Range("E8").Value = "  -1.17%  "

# Row 9
$ws.Range("D9").Value = "'0.3123"
$ws.Range("E9").Value = "  +3.10%  "

# Row 10
$ws.Range("D10").Value = "'25.19"
$ws.Range("E10").Value = "  +6.67%  "

# Row 11
$ws.Range("D11").Value = "'0.08286"
$ws.Range("E11").Value = "  +1.21%  "

# Row 12
$ws.Range("D12").Value = "'1.880.29"
$ws.Range("E12").Value = "  +2.13%  "

# Row 13
$ws.Range("D13").Value = "'0.7295"
$ws.Range("E13").Value = "  +3.15%  "

# Row 14
$ws.Range("D14").Value = "'5.299"
$ws.Range("E14").Value = "  +1.89%  "

# Row 15
$ws.Range("D15").Value = "'91.30"
$ws.Range("E15").Value = "  +1.68%  "

# Row 16
$ws.Range("D16").Value = "'29.463.64"
$ws.Range("E16").Value = "  +1.03%  "

# Row 17
$ws.Range("D17").Value = "'5.942"
$ws.Range("E17").Value = "  +1.95%  "

# Row 18
$ws.Range("D18").Value = "'247.67"
$ws.Range("E18").Value = "  +4.61%  "

# Row 19
$ws.Range("D19").Value = "'0.000007874"
$ws.Range("E19").Value = "  +0.22%  "

# Row 20
$ws.Range("D20").Value = "'13.36"
$ws.Range("E20").Value = "  +0.98%  "

# Row 21
$ws.Range("D21").Value = "'2.131.99"
$ws.Range("E21").Value = "  +2.62%  "

# Row 22
$ws.Range("D22").Value = "'0.9999"
$ws.Range("E22").Value = "  +0.02%  "

# Row 23
$ws.Range("D23").Value = "'7.991"
$ws.Range("E23").Value = "  +6.18%  "

# Row 24
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  -0.07%  "

# Row 25
$ws.Range("D25").Value = "'0.1599"
$ws.Range("E25").Value = "  +13.45%  "

# Row 26
$ws.Range("D26").Value = "'164.26"
$ws.Range("E26").Value = "  +0.61%  "

# Row 27
$ws.Range("D27").Value = "'9.031"
$ws.Range("E27").Value = "  +1.55%  "

# Row 28
$ws.Range("E28").Value = "  +1.33%  "

# Row 29
$ws.Range("D29").Value = "'1.362"
$ws.Range("E29").Value = "  -2.87%  "

# Row 30
$ws.Range("D30").Value = "'1.502"
$ws.Range("E30").Value = "  +2.05%  "

# Row 31
$ws.Range("D31").Value = "'4.408"
$ws.Range("E31").Value = "  +1.23%  "

# Row 32
$ws.Range("D32").Value = "'4.127"
$ws.Range("E32").Value = "  +2.41%  "

# Row 33
$ws.Range("D33").Value = "'0.05307"
$ws.Range("E33").Value = "  +2.24%  "

# Row 34
$ws.Range("D34").Value = "'1.943"
$ws.Range("E34").Value = "  +1.76%  "

# Row 35
$ws.Range("D35").Value = "'1.197"
$ws.Range("E35").Value = "  +2.67%  "

# Row 36
$ws.Range("D36").Value = "'0.7278"
$ws.Range("E36").Value = "  +1.75%  "

# Row 37
$ws.Range("D37").Value = "'2.679"
$ws.Range("E37").Value = "  -0.05%  "

# Row 38
$ws.Range("D38").Value = "'0.01872"
$ws.Range("E38").Value = "  +1.37%  "

# Row 39
$ws.Range("D39").Value = "'1.237.95"
$ws.Range("E39").Value = "  +7.48%  "

# Row 40
$ws.Range("D40").Value = "'2.725"
$ws.Range("E40").Value = "  +0.13%  "

# Row 41
$ws.Range("D41").Value = "'0.9134"
$ws.Range("E41").Value = "  -2.37%  "

# Row 42
$ws.Range("D42").Value = "'74.66"
$ws.Range("E42").Value = "  +6.32%  "

# Row 43
$ws.Range("D43").Value = "'6.196"
$ws.Range("E43").Value = "  +3.31%  "

# Row 44
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  -0.03%  "

# Row 45
$ws.Range("D45").Value = "'103.02"
$ws.Range("E45").Value = "  +0.07%  "

# Row 46
$ws.Range("D46").Value = "'2.031.58"
$ws.Range("E46").Value = "  +1.53%  "

# Row 47
$ws.Range("D47").Value = "'0.5300"
$ws.Range("E47").Value = "  +0.16%  "

# Row 48
$ws.Range("D48").Value = "'2.979"
$ws.Range("E48").Value = "  +14.91%  "

# Row 49
$ws.Range("D49").Value = "'1.770"
$ws.Range("E49").Value = "  +1.41%  "

# Row 50
$ws.Range("D50").Value = "'9.323"
$ws.Range("E50").Value = "  +1.86%  "

# Row 51 - coin changed from TheSandbox to BabyDogeCoin
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "'0.00000000119"
$ws.Range("E51").Value = "  +0.25%  "
